$d = $word.ActiveDocument
$startPara = $d.Paragraphs.Item(2)
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng = $d.Range($startPara.Range.Start, $lastPara.Range.End)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Heading4"/>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:t>IWATER</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>I</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve">WATER </w:t>
  </w:r>
  <w:r>
    <w:t>is the interface class containing the virtual functions relating which is implemented in the hardware class.</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
</w:p>
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="GridTable1Light"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="4675"/>
    <w:gridCol w:w="4675"/>
  </w:tblGrid>
  <w:tr>
    <w:trPr>
      <w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="9350" w:type="dxa"/>
        <w:gridSpan w:val="2"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t>IWATER</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="4675" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="D0CECE" w:themeFill="background2" w:themeFillShade="E6"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:b w:val="0"/>
            <w:bCs w:val="0"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b w:val="0"/>
            <w:bCs w:val="0"/>
          </w:rPr>
          <w:t>Functions</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="D0CECE" w:themeFill="background2" w:themeFillShade="E6"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="both"/>
          <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t xml:space="preserve">Description </w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:i/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t>virtual voi</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t>d SetWaterLevel(int level) = 0;</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>The purpose of this function will be to</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> set the appropriate water lever for the washing process</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve">. </w:t>
        </w:r>
        <w:r>
          <w:t>It takes one parameter for the level that we want to check.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:i/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t>virtual boolean CheckWaterLevel(int level) = 0;</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:i/>
          </w:rPr>
        </w:pPr>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>The purpose of this function will be to</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> check if the washing machine has reached the desired water level</w:t>
        </w:r>
        <w:r>
          <w:t>. It takes one parameter for the level that we want to check.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:i/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t>virtual void SinkWater() = 0;</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t xml:space="preserve">The purpose of this function will be to </w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="9350" w:type="dxa"/>
        <w:gridSpan w:val="2"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t xml:space="preserve">Remarks: </w:t>
        </w:r>
        <w:r>
          <w:t>This interface is complete. All the func</w:t>
        </w:r>
        <w:r>
          <w:t>t</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve">ions are implemented. </w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t>To discuss: I just wonder if we need any more functions for this interface or are those 3 enough?</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t>(Please describe here any problems, remarks, current states of your functions)</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:br w:type="page"/>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="Heading4"/>
    <w:jc w:val="both"/>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>IMOTOR</w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:t>I</w:t>
  </w:r>
  <w:r>
    <w:t>MOTOR</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> is the interface class containing the virtual functions relating which is implemented in the hardware class.</w:t>
  </w:r>
</w:p>
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:tblPr>
    <w:tblStyle w:val="GridTable1Light"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="4675"/>
    <w:gridCol w:w="4675"/>
  </w:tblGrid>
  <w:tr>
    <w:trPr>
      <w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
    </w:trPr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="9350" w:type="dxa"/>
        <w:gridSpan w:val="2"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t>IMOTOR</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="4675" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="D0CECE" w:themeFill="background2" w:themeFillShade="E6"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:b w:val="0"/>
            <w:bCs w:val="0"/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b w:val="0"/>
            <w:bCs w:val="0"/>
          </w:rPr>
          <w:t>Functions</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
        <w:shd w:val="clear" w:color="auto" w:fill="D0CECE" w:themeFill="background2" w:themeFillShade="E6"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="both"/>
          <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t xml:space="preserve">Description </w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:i/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t>virtual v</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t>oid SetDirection(char dir) = 0;</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t xml:space="preserve">The purpose of this function will be to </w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve">indicate the direction of the washing machine. </w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve">It takes one parameter </w:t>
        </w:r>
        <w:r>
          <w:t>to allow us to indicate t</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve">he </w:t>
        </w:r>
        <w:r>
          <w:t>direction</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> that we want to </w:t>
        </w:r>
        <w:r>
          <w:t>set</w:t>
        </w:r>
        <w:r>
          <w:t>.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:i/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t>virtual</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t xml:space="preserve"> void SetSpeed(int level</w:t>
        </w:r>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t xml:space="preserve">) = 0; </w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t>The purp</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve">ose of this function </w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve">will be to indicate the </w:t>
        </w:r>
        <w:r>
          <w:t>speed</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> of the washing machine. It takes one parameter to allow us to indicate the </w:t>
        </w:r>
        <w:r>
          <w:t>speed</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> that we want to set.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:i/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t>virtual void StopMotor() = 0;</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t xml:space="preserve">The purpose of this function will be to </w:t>
        </w:r>
        <w:r>
          <w:t>stop the motor.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:i/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t>virtual void StartMotor() = 0</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t xml:space="preserve">The purpose of this function will be to </w:t>
        </w:r>
        <w:r>
          <w:t>start the motor</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> at an initial state assuming that all the prerequisites are handled properly</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> within the IProgram interface</w:t>
        </w:r>
        <w:bookmarkStart w:id="0" w:name="_GoBack"/>
        <w:bookmarkEnd w:id="0"/>
        <w:r>
          <w:t>.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
          <w:jc w:val="both"/>
          <w:rPr>
            <w:i/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:i/>
          </w:rPr>
          <w:t>virtual void CheckLoadingLevel(int level) = 0;</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="4675" w:type="dxa"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:cnfStyle w:val="000000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        </w:pPr>
        <w:r>
          <w:t xml:space="preserve">The purpose of this function will be to </w:t>
        </w:r>
        <w:r>
          <w:t>check le loading level of the washing machine. It takes one parameter for the level that we want to check.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/>
        <w:tcW w:w="9350" w:type="dxa"/>
        <w:gridSpan w:val="2"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t xml:space="preserve">Remarks: </w:t>
        </w:r>
        <w:r>
          <w:t>This interface is complete. All the func</w:t>
        </w:r>
        <w:r>
          <w:t>t</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve">ions are implemented. </w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t xml:space="preserve">To discuss: </w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve">To start the motor can we just set the initial speed, or we also </w:t>
        </w:r>
        <w:r>
          <w:t>should</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> consider some othe</w:t>
        </w:r>
        <w:r>
          <w:t>r</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> states?</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> Should we take care of the </w:t>
        </w:r>
        <w:r>
          <w:t>prerequisites</w:t>
        </w:r>
        <w:r>
          <w:t xml:space="preserve"> in this function of in the IProgram interface?</w:t>
        </w:r>
      </w:p>
      <w:p>
        <w:pPr>
          <w:jc w:val="both"/>
        </w:pPr>
        <w:r>
          <w:t>(Please describe here any problems, remarks, current states of your functions)</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@

$rng.InsertXML($xml)
